# Add Week 15 simulations: two new WR rows (D.Davis, P.Dorsett) with zeroed stats,
# and update the active sheet/selection state to reflect the WR sheet being the
# last one worked on.

$wb = $excel.ActiveWorkbook

$wsRB = $wb.Worksheets.Item("RB")
$wsWR = $wb.Worksheets.Item("WR")

# New players added to the WR sheet (rows 9 and 10)
$wsWR.Range("A9").Value = "D.Davis"
$wsWR.Range("B9").Value = 0
$wsWR.Range("C9").Value = 0
$wsWR.Range("D9").Value = 0
$wsWR.Range("E9").Value = 0
$wsWR.Range("F9").Value = 0
$wsWR.Range("G9").Value = 0
$wsWR.Range("H9").Value = 0
$wsWR.Range("I9").Value = 0
$wsWR.Range("J9").Value = 0

$wsWR.Range("A10").Value = "P.Dorsett"
$wsWR.Range("B10").Value = 0
$wsWR.Range("C10").Value = 0
$wsWR.Range("D10").Value = 0
$wsWR.Range("E10").Value = 0
$wsWR.Range("F10").Value = 0
$wsWR.Range("G10").Value = 0
$wsWR.Range("H10").Value = 0
$wsWR.Range("I10").Value = 0
$wsWR.Range("J10").Value = 0

# Selection on RB was left at B25 before switching to WR
$wsRB.Activate()
$wsRB.Range("B25").Select() | Out-Null

# WR ends up as the active (selected) tab, with J11 selected
$wsWR.Activate()
$wsWR.Range("J11").Select() | Out-Null
